# Update the "想去人数" (interested count) figures in the F column
# for both the "展览" and "全部类型" sheets (which mirror each other).
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 1745
    $ws.Range("F3").Value = 8000
    $ws.Range("F4").Value = 185
    $ws.Range("F5").Value = 278
}
